$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 41129
$ws.Range("E2").Value = 11719
$ws.Range("F2").Value = 11719
$ws.Range("G2").Value = 11579
$ws.Range("H2").Value = 8138
$ws.Range("I2").Value = 8256
$ws.Range("J2").Value = -118
$ws.Range("K2").Value = 74187
$ws.Range("L2").Value = 17106
$ws.Range("M2").Value = 57081
$ws.Range("N2").Value = 56306
$ws.Range("O2").Value = 775
$ws.Range("P2").Value = 9550
$ws.Range("Q2").Value = 7652
$ws.Range("R2").Value = -3429
$ws.Range("S2").Value = -3538
$ws.Range("T2").Value = 3304
$ws.Range("U2").Value = 4348
$ws.Range("V2").Value = 2809
$ws.Range("W2").Value = 28.49
$ws.Range("X2").Value = 19.79
$ws.Range("Y2").Value = 15.14
$ws.Range("Z2").Value = 11.23
$ws.Range("AA2").Value = 29.97
$ws.Range("AB2").Value = 535.0700000000001
$ws.Range("AC2").Value = 6013
$ws.Range("AD2").Value = 12.65
$ws.Range("AE2").Value = 44723
$ws.Range("AF2").Value = 1.7
$ws.Range("AG2").Value = 3400
$ws.Range("AH2").Value = 4.47
$ws.Range("AI2").Value = 51.85
$ws.Range("AJ2").Value = 137292497

# Row 3
$ws.Range("D3").Value = 41698
$ws.Range("E3").Value = 13659
$ws.Range("F3").Value = 13659
$ws.Range("G3").Value = 14144
$ws.Range("H3").Value = 10322
$ws.Range("I3").Value = 10357
$ws.Range("J3").Value = -35
$ws.Range("K3").Value = 86734
$ws.Range("L3").Value = 23979
$ws.Range("M3").Value = 62755
$ws.Range("N3").Value = 62089
$ws.Range("O3").Value = 666
$ws.Range("P3").Value = 9550
$ws.Range("Q3").Value = 12592
$ws.Range("R3").Value = -7476
$ws.Range("S3").Value = -3810
$ws.Range("T3").Value = 2102
$ws.Range("U3").Value = 10490
$ws.Range("V3").Value = 3363
$ws.Range("W3").Value = 32.76
$ws.Range("X3").Value = 24.75
$ws.Range("Y3").Value = 17.5
$ws.Range("Z3").Value = 12.83
$ws.Range("AA3").Value = 38.21
$ws.Range("AB3").Value = 598.61
$ws.Range("AC3").Value = 7544
$ws.Range("AD3").Value = 13.85
$ws.Range("AE3").Value = 49290
$ws.Range("AF3").Value = 2.12
$ws.Range("AG3").Value = 3400
$ws.Range("AH3").Value = 3.25
$ws.Range("AI3").Value = 41.35
$ws.Range("AJ3").Value = 137292497

# Row 4
$ws.Range("D4").Value = 45033
$ws.Range("E4").Value = 14696
$ws.Range("F4").Value = 14701
$ws.Range("G4").Value = 15873
$ws.Range("H4").Value = 12255
$ws.Range("I4").Value = 12308
$ws.Range("J4").Value = -52
$ws.Range("K4").Value = 98121
$ws.Range("L4").Value = 26943
$ws.Range("M4").Value = 71178
$ws.Range("N4").Value = 70452
$ws.Range("O4").Value = 726
$ws.Range("P4").Value = 9550
$ws.Range("Q4").Value = 14978
$ws.Range("R4").Value = -7026
$ws.Range("S4").Value = -4884
$ws.Range("T4").Value = 1617
$ws.Range("U4").Value = 13360
$ws.Range("V4").Value = 2621
$ws.Range("W4").Value = 32.63
$ws.Range("X4").Value = 27.21
$ws.Range("Y4").Value = 18.57
$ws.Range("Z4").Value = 13.26
$ws.Range("AA4").Value = 37.85
$ws.Range("AB4").Value = 682.6799999999999
$ws.Range("AC4").Value = 8965
$ws.Range("AD4").Value = 11.27
$ws.Range("AE4").Value = 55797
$ws.Range("AF4").Value = 1.81
$ws.Range("AG4").Value = 3600
$ws.Range("AH4").Value = 3.56
$ws.Range("AI4").Value = 36.93
$ws.Range("AJ4").Value = 137292497

# Row 5
$ws.Range("D5").Value = 46672
$ws.Range("E5").Value = 14261
$ws.Range("F5").Value = 14261
$ws.Range("G5").Value = 13456
$ws.Range("H5").Value = 11642
$ws.Range("I5").Value = 11638
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 99814
$ws.Range("L5").Value = 21516
$ws.Range("M5").Value = 78298
$ws.Range("N5").Value = 77721
$ws.Range("O5").Value = 577
$ws.Range("P5").Value = 9550
$ws.Range("Q5").Value = 11431
$ws.Range("R5").Value = -6830
$ws.Range("S5").Value = -4386
$ws.Range("T5").Value = 3183
$ws.Range("U5").Value = 8248
$ws.Range("V5").Value = 2789
$ws.Range("W5").Value = 30.56
$ws.Range("X5").Value = 24.95
$ws.Range("Y5").Value = 15.71
$ws.Range("Z5").Value = 11.76
$ws.Range("AA5").Value = 27.48
$ws.Range("AB5").Value = 757.02
$ws.Range("AC5").Value = 8476
$ws.Range("AD5").Value = 13.63
$ws.Range("AE5").Value = 61553
$ws.Range("AF5").Value = 1.88
$ws.Range("AG5").Value = 4000
$ws.Range("AH5").Value = 3.46
$ws.Range("AI5").Value = 43.4
$ws.Range("AJ5").Value = 137292497

# Row 6
$ws.Range("D6").Value = 44715
$ws.Range("E6").Value = 12551
$ws.Range("F6").Value = 12551
$ws.Range("G6").Value = 13187
$ws.Range("H6").Value = 8987
$ws.Range("I6").Value = 9016
$ws.Range("K6").Value = 101551
$ws.Range("L6").Value = 20220
$ws.Range("M6").Value = 81331
$ws.Range("N6").Value = 80793
$ws.Range("P6").Value = 9550
$ws.Range("Q6").Value = 8221
$ws.Range("R6").Value = -463
$ws.Range("S6").Value = -5498
$ws.Range("T6").Value = 3939
$ws.Range("U6").Value = 4282
$ws.Range("V6").Value = 2284
$ws.Range("W6").Value = 28.07
$ws.Range("X6").Value = 20.1
$ws.Range("Y6").Value = 11.38
$ws.Range("Z6").Value = 8.93
$ws.Range("AA6").Value = 24.86
$ws.Range("AB6").Value = 795.11
$ws.Range("AC6").Value = 6567
$ws.Range("AD6").Value = 15.46
$ws.Range("AE6").Value = 63987
$ws.Range("AF6").Value = 1.59
$ws.Range("AG6").Value = 4000
$ws.Range("AH6").Value = 3.94
$ws.Range("AI6").Value = 56.02
$ws.Range("AJ6").Value = 137292497

# Row 7
$ws.Range("D7").Value = 49608
$ws.Range("E7").Value = 14306
$ws.Range("G7").Value = 15914
$ws.Range("H7").Value = 11371
$ws.Range("I7").Value = 11347
$ws.Range("K7").Value = 108515
$ws.Range("L7").Value = 21153
$ws.Range("M7").Value = 87362
$ws.Range("N7").Value = 86787
$ws.Range("P7").Value = 9550
$ws.Range("Q7").Value = 10623
$ws.Range("R7").Value = -4228
$ws.Range("S7").Value = -5451
$ws.Range("T7").Value = 3442
$ws.Range("U7").Value = 6021
$ws.Range("W7").Value = 28.84
$ws.Range("X7").Value = 22.92
$ws.Range("Y7").Value = 13.54
$ws.Range("Z7").Value = 10.83
$ws.Range("AA7").Value = 24.21
$ws.Range("AC7").Value = 8265
$ws.Range("AD7").Value = 11.51
$ws.Range("AE7").Value = 68563
$ws.Range("AF7").Value = 1.39
$ws.Range("AG7").Value = 4350
$ws.Range("AH7").Value = 4.57
$ws.Range("AI7").Value = 52.63

# Row 8
$ws.Range("D8").Value = 53675
$ws.Range("E8").Value = 15642
$ws.Range("G8").Value = 16521
$ws.Range("H8").Value = 11784
$ws.Range("I8").Value = 11756
$ws.Range("K8").Value = 114889
$ws.Range("L8").Value = 21906
$ws.Range("M8").Value = 92983
$ws.Range("N8").Value = 92413
$ws.Range("P8").Value = 9550
$ws.Range("Q8").Value = 10143
$ws.Range("R8").Value = -4542
$ws.Range("S8").Value = -5712
$ws.Range("T8").Value = 3668
$ws.Range("U8").Value = 6748
$ws.Range("W8").Value = 29.14
$ws.Range("X8").Value = 21.95
$ws.Range("Y8").Value = 13.12
$ws.Range("Z8").Value = 10.55
$ws.Range("AA8").Value = 23.56
$ws.Range("AC8").Value = 8563
$ws.Range("AD8").Value = 11.11
$ws.Range("AE8").Value = 73008
$ws.Range("AF8").Value = 1.3
$ws.Range("AG8").Value = 4479
$ws.Range("AH8").Value = 4.71
$ws.Range("AI8").Value = 52.3

# Row 9
$ws.Range("D9").Value = 56488
$ws.Range("E9").Value = 16327
$ws.Range("G9").Value = 17242
$ws.Range("H9").Value = 12281
$ws.Range("I9").Value = 12254
$ws.Range("K9").Value = 121278
$ws.Range("L9").Value = 22462
$ws.Range("M9").Value = 98815
$ws.Range("N9").Value = 98300
$ws.Range("P9").Value = 9550
$ws.Range("Q9").Value = 11502
$ws.Range("R9").Value = -4102
$ws.Range("S9").Value = -5358
$ws.Range("T9").Value = 3684
$ws.Range("U9").Value = 7819
$ws.Range("W9").Value = 28.9
$ws.Range("X9").Value = 21.74
$ws.Range("Y9").Value = 12.85
$ws.Range("Z9").Value = 10.4
$ws.Range("AA9").Value = 22.73
$ws.Range("AC9").Value = 8925
$ws.Range("AD9").Value = 10.66
$ws.Range("AE9").Value = 77659
$ws.Range("AF9").Value = 1.22
$ws.Range("AG9").Value = 4636
$ws.Range("AH9").Value = 4.87
$ws.Range("AI9").Value = 51.94
